$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '25.769.20'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  -0.90%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.627.33'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  -0.79%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.004'
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  +0.29%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '215.38'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  +0.11%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.5055'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  -0.82%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '1.004'
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  +0.28%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.2569'
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  -0.51%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.06417'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  +0.89%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '19.35'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  -2.76%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07782'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  +0.44%  '
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  -1.08%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '1.629.83'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  -0.58%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '1.851.97'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  -0.84%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.5568'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  +1.66%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '62.91'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  -2.31%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.0₅7542'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  -2.73%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '25.792.21'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  -0.90%  '
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  +0.18%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '193.50'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  -1.83%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '4.315'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  -3.34%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '9.812'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  -1.56%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '5.989'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  -2.42%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '1.004'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  +0.20%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '1.799'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  -4.93%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '140.86'
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  -1.38%  '
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  -0.14%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '6.724'
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  -2.24%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '15.37'
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  -1.62%  '
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  -0.33%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.04864'
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  -0.59%  '
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  -0.30%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '3.189'
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  -0.81%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.554'
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  -0.22%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '2.377'
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  +0.10%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.8928'
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  -2.89%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.564'
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  +0.19%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '1.127.15'
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  +2.04%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.5460'
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  -1.83%  '
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  -0.76%  '
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  +0.02%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '5.552'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  -1.13%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.7962'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  -1.06%  '
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  -1.54%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '1.779.67'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  -0.02%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.0₈113'
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  -6.93%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.4440'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  -2.07%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '55.15'
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  -0.41%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.05053'
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  -2.60%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '7.633'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  +0.43%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.001'
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  -0.11%  '
